$wb = $excel.ActiveWorkbook
$dayWs = $wb.Worksheets.Item("day")
$weekWs = $wb.Worksheets.Item("week")

# --- Append 25 new rows (584-608) to the 'day' sheet ---

# Row 584: MRF
$dayWs.Cells.Item(584,1).Value = 1
$dayWs.Cells.Item(584,2).Value = 'MRF'
$dayWs.Cells.Item(584,3).Value = 'Mrf Limited'
$dayWs.Cells.Item(584,4).NumberFormat = "@"
$dayWs.Cells.Item(584,4).Value = '500290'
$dayWs.Cells.Item(584,4).ClearFormats()
$dayWs.Cells.Item(584,5).Value = 1
$dayWs.Cells.Item(584,6).Value = 137292.25
$dayWs.Cells.Item(584,7).Value = 5664
$dayWs.Cells.Item(584,8).Value = 'day'
$dayWs.Cells.Item(584,9).Value = '23/09/2024 11:34:50'

# Row 585: BANKNIFTY
$dayWs.Cells.Item(585,1).Value = 2
$dayWs.Cells.Item(585,2).Value = 'BANKNIFTY'
$dayWs.Cells.Item(585,3).Value = 'BANKNIFTY'
# D585 left empty (matches source data which has no bsecode for this index)
$dayWs.Cells.Item(585,5).Value = 0.58
$dayWs.Cells.Item(585,6).Value = 54105.8
$dayWs.Cells.Item(585,7).Value = 0
$dayWs.Cells.Item(585,8).Value = 'day'
$dayWs.Cells.Item(585,9).Value = '23/09/2024 11:34:50'

# Row 586: PAGEIND
$dayWs.Cells.Item(586,1).Value = 3
$dayWs.Cells.Item(586,2).Value = 'PAGEIND'
$dayWs.Cells.Item(586,3).Value = 'Page Industries Limited'
$dayWs.Cells.Item(586,4).NumberFormat = "@"
$dayWs.Cells.Item(586,4).Value = '532827'
$dayWs.Cells.Item(586,4).ClearFormats()
$dayWs.Cells.Item(586,5).Value = -0.08
$dayWs.Cells.Item(586,6).Value = 43017.65
$dayWs.Cells.Item(586,7).Value = 47601
$dayWs.Cells.Item(586,8).Value = 'day'
$dayWs.Cells.Item(586,9).Value = '23/09/2024 11:34:50'

# Row 587: NIFTY
$dayWs.Cells.Item(587,1).Value = 4
$dayWs.Cells.Item(587,2).Value = 'NIFTY'
$dayWs.Cells.Item(587,3).Value = 'NIFTY'
# D587 left empty (matches source data which has no bsecode for this index)
$dayWs.Cells.Item(587,5).Value = 0.57
$dayWs.Cells.Item(587,6).Value = 25939.05
$dayWs.Cells.Item(587,7).Value = 0
$dayWs.Cells.Item(587,8).Value = 'day'
$dayWs.Cells.Item(587,9).Value = '23/09/2024 11:34:50'

# Row 588: DIXON
$dayWs.Cells.Item(588,1).Value = 5
$dayWs.Cells.Item(588,2).Value = 'DIXON'
$dayWs.Cells.Item(588,3).Value = 'Dixon Technologies'
$dayWs.Cells.Item(588,4).NumberFormat = "@"
$dayWs.Cells.Item(588,4).Value = '540699'
$dayWs.Cells.Item(588,4).ClearFormats()
$dayWs.Cells.Item(588,5).Value = 1.85
$dayWs.Cells.Item(588,6).Value = 14254.85
$dayWs.Cells.Item(588,7).Value = 464206
$dayWs.Cells.Item(588,8).Value = 'day'
$dayWs.Cells.Item(588,9).Value = '23/09/2024 11:34:50'

# Row 589: SIEMENS
$dayWs.Cells.Item(589,1).Value = 6
$dayWs.Cells.Item(589,2).Value = 'SIEMENS'
$dayWs.Cells.Item(589,3).Value = 'Siemens Limited'
$dayWs.Cells.Item(589,4).NumberFormat = "@"
$dayWs.Cells.Item(589,4).Value = '500550'
$dayWs.Cells.Item(589,4).ClearFormats()
$dayWs.Cells.Item(589,5).Value = 1.55
$dayWs.Cells.Item(589,6).Value = 6938.4
$dayWs.Cells.Item(589,7).Value = 215612
$dayWs.Cells.Item(589,8).Value = 'day'
$dayWs.Cells.Item(589,9).Value = '23/09/2024 11:34:50'

# Row 590: ESCORTS
$dayWs.Cells.Item(590,1).Value = 7
$dayWs.Cells.Item(590,2).Value = 'ESCORTS'
$dayWs.Cells.Item(590,3).Value = 'Escorts Limited'
$dayWs.Cells.Item(590,4).NumberFormat = "@"
$dayWs.Cells.Item(590,4).Value = '500495'
$dayWs.Cells.Item(590,4).ClearFormats()
$dayWs.Cells.Item(590,5).Value = 3.79
$dayWs.Cells.Item(590,6).Value = 4233.9
$dayWs.Cells.Item(590,7).Value = 543509
$dayWs.Cells.Item(590,8).Value = 'day'
$dayWs.Cells.Item(590,9).Value = '23/09/2024 11:34:50'

# Row 591: TVSMOTOR
$dayWs.Cells.Item(591,1).Value = 8
$dayWs.Cells.Item(591,2).Value = 'TVSMOTOR'
$dayWs.Cells.Item(591,3).Value = 'Tvs Motor Company Limited'
$dayWs.Cells.Item(591,4).NumberFormat = "@"
$dayWs.Cells.Item(591,4).Value = '532343'
$dayWs.Cells.Item(591,4).ClearFormats()
$dayWs.Cells.Item(591,5).Value = 1.05
$dayWs.Cells.Item(591,6).Value = 2845.2
$dayWs.Cells.Item(591,7).Value = 1352766
$dayWs.Cells.Item(591,8).Value = 'day'
$dayWs.Cells.Item(591,9).Value = '23/09/2024 11:34:50'

# Row 592: HAVELLS
$dayWs.Cells.Item(592,1).Value = 9
$dayWs.Cells.Item(592,2).Value = 'HAVELLS'
$dayWs.Cells.Item(592,3).Value = 'Havells India Limited'
$dayWs.Cells.Item(592,4).NumberFormat = "@"
$dayWs.Cells.Item(592,4).Value = '517354'
$dayWs.Cells.Item(592,4).ClearFormats()
$dayWs.Cells.Item(592,5).Value = 1.67
$dayWs.Cells.Item(592,6).Value = 2082.4
$dayWs.Cells.Item(592,7).Value = 1311949
$dayWs.Cells.Item(592,8).Value = 'day'
$dayWs.Cells.Item(592,9).Value = '23/09/2024 11:34:50'

# Row 593: TECHM
$dayWs.Cells.Item(593,1).Value = 10
$dayWs.Cells.Item(593,2).Value = 'TECHM'
$dayWs.Cells.Item(593,3).Value = 'Tech Mahindra Limited'
$dayWs.Cells.Item(593,4).NumberFormat = "@"
$dayWs.Cells.Item(593,4).Value = '532755'
$dayWs.Cells.Item(593,4).ClearFormats()
$dayWs.Cells.Item(593,5).Value = -0.92
$dayWs.Cells.Item(593,6).Value = 1607.15
$dayWs.Cells.Item(593,7).Value = 1535903
$dayWs.Cells.Item(593,8).Value = 'day'
$dayWs.Cells.Item(593,9).Value = '23/09/2024 11:34:50'

# Row 594: BHARATFORG
$dayWs.Cells.Item(594,1).Value = 11
$dayWs.Cells.Item(594,2).Value = 'BHARATFORG'
$dayWs.Cells.Item(594,3).Value = 'Bharat Forge Limited'
$dayWs.Cells.Item(594,4).NumberFormat = "@"
$dayWs.Cells.Item(594,4).Value = '500493'
$dayWs.Cells.Item(594,4).ClearFormats()
$dayWs.Cells.Item(594,5).Value = -0.58
$dayWs.Cells.Item(594,6).Value = 1582.2
$dayWs.Cells.Item(594,7).Value = 2341523
$dayWs.Cells.Item(594,8).Value = 'day'
$dayWs.Cells.Item(594,9).Value = '23/09/2024 11:34:50'

# Row 595: BATAINDIA
$dayWs.Cells.Item(595,1).Value = 12
$dayWs.Cells.Item(595,2).Value = 'BATAINDIA'
$dayWs.Cells.Item(595,3).Value = 'Bata India Limited'
$dayWs.Cells.Item(595,4).NumberFormat = "@"
$dayWs.Cells.Item(595,4).Value = '500043'
$dayWs.Cells.Item(595,4).ClearFormats()
$dayWs.Cells.Item(595,5).Value = 0.19
$dayWs.Cells.Item(595,6).Value = 1428.8
$dayWs.Cells.Item(595,7).Value = 112972
$dayWs.Cells.Item(595,8).Value = 'day'
$dayWs.Cells.Item(595,9).Value = '23/09/2024 11:34:50'

# Row 596: ICICIBANK
$dayWs.Cells.Item(596,1).Value = 13
$dayWs.Cells.Item(596,2).Value = 'ICICIBANK'
$dayWs.Cells.Item(596,3).Value = 'Icici Bank Limited'
$dayWs.Cells.Item(596,4).NumberFormat = "@"
$dayWs.Cells.Item(596,4).Value = '532174'
$dayWs.Cells.Item(596,4).ClearFormats()
$dayWs.Cells.Item(596,5).Value = -1.24
$dayWs.Cells.Item(596,6).Value = 1321.9
$dayWs.Cells.Item(596,7).Value = 11380409
$dayWs.Cells.Item(596,8).Value = 'day'
$dayWs.Cells.Item(596,9).Value = '23/09/2024 11:34:50'

# Row 597: AXISBANK
$dayWs.Cells.Item(597,1).Value = 14
$dayWs.Cells.Item(597,2).Value = 'AXISBANK'
$dayWs.Cells.Item(597,3).Value = 'Axis Bank Limited'
$dayWs.Cells.Item(597,4).NumberFormat = "@"
$dayWs.Cells.Item(597,4).Value = '532215'
$dayWs.Cells.Item(597,4).ClearFormats()
$dayWs.Cells.Item(597,5).Value = 0.14
$dayWs.Cells.Item(597,6).Value = 1246.8
$dayWs.Cells.Item(597,7).Value = 5462814
$dayWs.Cells.Item(597,8).Value = 'day'
$dayWs.Cells.Item(597,9).Value = '23/09/2024 11:34:50'

# Row 598: JUBLFOOD
$dayWs.Cells.Item(598,1).Value = 15
$dayWs.Cells.Item(598,2).Value = 'JUBLFOOD'
$dayWs.Cells.Item(598,3).Value = 'Jubilant Foodworks Limited'
$dayWs.Cells.Item(598,4).NumberFormat = "@"
$dayWs.Cells.Item(598,4).Value = '533155'
$dayWs.Cells.Item(598,4).ClearFormats()
$dayWs.Cells.Item(598,5).Value = 0.21
$dayWs.Cells.Item(598,6).Value = 704.85
$dayWs.Cells.Item(598,7).Value = 1974227
$dayWs.Cells.Item(598,8).Value = 'day'
$dayWs.Cells.Item(598,9).Value = '23/09/2024 11:34:50'

# Row 599: TATAPOWER
$dayWs.Cells.Item(599,1).Value = 16
$dayWs.Cells.Item(599,2).Value = 'TATAPOWER'
$dayWs.Cells.Item(599,3).Value = 'Tata Power Company Limited'
$dayWs.Cells.Item(599,4).NumberFormat = "@"
$dayWs.Cells.Item(599,4).Value = '500400'
$dayWs.Cells.Item(599,4).ClearFormats()
$dayWs.Cells.Item(599,5).Value = 2.32
$dayWs.Cells.Item(599,6).Value = 454.45
$dayWs.Cells.Item(599,7).Value = 26480471
$dayWs.Cells.Item(599,8).Value = 'day'
$dayWs.Cells.Item(599,9).Value = '23/09/2024 11:34:50'

# Row 600: NTPC
$dayWs.Cells.Item(600,1).Value = 17
$dayWs.Cells.Item(600,2).Value = 'NTPC'
$dayWs.Cells.Item(600,3).Value = 'Ntpc Limited'
$dayWs.Cells.Item(600,4).NumberFormat = "@"
$dayWs.Cells.Item(600,4).Value = '532555'
$dayWs.Cells.Item(600,4).ClearFormats()
$dayWs.Cells.Item(600,5).Value = 1.04
$dayWs.Cells.Item(600,6).Value = 428.35
$dayWs.Cells.Item(600,7).Value = 14139629
$dayWs.Cells.Item(600,8).Value = 'day'
$dayWs.Cells.Item(600,9).Value = '23/09/2024 11:34:50'

# Row 601: ABFRL
$dayWs.Cells.Item(601,1).Value = 18
$dayWs.Cells.Item(601,2).Value = 'ABFRL'
$dayWs.Cells.Item(601,3).Value = 'Aditya Birla Fashion And Retail Limited'
$dayWs.Cells.Item(601,4).NumberFormat = "@"
$dayWs.Cells.Item(601,4).Value = '535755'
$dayWs.Cells.Item(601,4).ClearFormats()
$dayWs.Cells.Item(601,5).Value = 5.11
$dayWs.Cells.Item(601,6).Value = 344.5
$dayWs.Cells.Item(601,7).Value = 12377319
$dayWs.Cells.Item(601,8).Value = 'day'
$dayWs.Cells.Item(601,9).Value = '23/09/2024 11:34:50'

# Row 602: POWERGRID
$dayWs.Cells.Item(602,1).Value = 19
$dayWs.Cells.Item(602,2).Value = 'POWERGRID'
$dayWs.Cells.Item(602,3).Value = 'Power Grid Corporation Of India Limited'
$dayWs.Cells.Item(602,4).NumberFormat = "@"
$dayWs.Cells.Item(602,4).Value = '532898'
$dayWs.Cells.Item(602,4).ClearFormats()
$dayWs.Cells.Item(602,5).Value = 0.09
$dayWs.Cells.Item(602,6).Value = 341.15
$dayWs.Cells.Item(602,7).Value = 6548154
$dayWs.Cells.Item(602,8).Value = 'day'
$dayWs.Cells.Item(602,9).Value = '23/09/2024 11:34:50'

# Row 603: HINDCOPPER
$dayWs.Cells.Item(603,1).Value = 20
$dayWs.Cells.Item(603,2).Value = 'HINDCOPPER'
$dayWs.Cells.Item(603,3).Value = 'Hindustan Copper Limited'
$dayWs.Cells.Item(603,4).NumberFormat = "@"
$dayWs.Cells.Item(603,4).Value = '513599'
$dayWs.Cells.Item(603,4).ClearFormats()
$dayWs.Cells.Item(603,5).Value = 0.27
$dayWs.Cells.Item(603,6).Value = 330.7
$dayWs.Cells.Item(603,7).Value = 4787855
$dayWs.Cells.Item(603,8).Value = 'day'
$dayWs.Cells.Item(603,9).Value = '23/09/2024 11:34:50'

# Row 604: NMDC
$dayWs.Cells.Item(604,1).Value = 21
$dayWs.Cells.Item(604,2).Value = 'NMDC'
$dayWs.Cells.Item(604,3).Value = 'Nmdc Limited'
$dayWs.Cells.Item(604,4).NumberFormat = "@"
$dayWs.Cells.Item(604,4).Value = '526371'
$dayWs.Cells.Item(604,4).ClearFormats()
$dayWs.Cells.Item(604,5).Value = 1.34
$dayWs.Cells.Item(604,6).Value = 215.39
$dayWs.Cells.Item(604,7).Value = 5902902
$dayWs.Cells.Item(604,8).Value = 'day'
$dayWs.Cells.Item(604,9).Value = '23/09/2024 11:34:50'

# Row 605: LTF
$dayWs.Cells.Item(605,1).Value = 22
$dayWs.Cells.Item(605,2).Value = 'LTF'
$dayWs.Cells.Item(605,3).Value = 'L&T Finance Ltd'
$dayWs.Cells.Item(605,4).NumberFormat = "@"
$dayWs.Cells.Item(605,4).Value = '533519'
$dayWs.Cells.Item(605,4).ClearFormats()
$dayWs.Cells.Item(605,5).Value = 1.77
$dayWs.Cells.Item(605,6).Value = 184.52
$dayWs.Cells.Item(605,7).Value = 13736132
$dayWs.Cells.Item(605,8).Value = 'day'
$dayWs.Cells.Item(605,9).Value = '23/09/2024 11:34:50'

# Row 606: NATIONALUM
$dayWs.Cells.Item(606,1).Value = 23
$dayWs.Cells.Item(606,2).Value = 'NATIONALUM'
$dayWs.Cells.Item(606,3).Value = 'National Aluminium Company Limited'
$dayWs.Cells.Item(606,4).NumberFormat = "@"
$dayWs.Cells.Item(606,4).Value = '532234'
$dayWs.Cells.Item(606,4).ClearFormats()
$dayWs.Cells.Item(606,5).Value = -2.06
$dayWs.Cells.Item(606,6).Value = 180.24
$dayWs.Cells.Item(606,7).Value = 13318424
$dayWs.Cells.Item(606,8).Value = 'day'
$dayWs.Cells.Item(606,9).Value = '23/09/2024 11:34:50'

# Row 607: CANBK
$dayWs.Cells.Item(607,1).Value = 24
$dayWs.Cells.Item(607,2).Value = 'CANBK'
$dayWs.Cells.Item(607,3).Value = 'Canara Bank'
$dayWs.Cells.Item(607,4).NumberFormat = "@"
$dayWs.Cells.Item(607,4).Value = '532483'
$dayWs.Cells.Item(607,4).ClearFormats()
$dayWs.Cells.Item(607,5).Value = 4.13
$dayWs.Cells.Item(607,6).Value = 109.3
$dayWs.Cells.Item(607,7).Value = 44004939
$dayWs.Cells.Item(607,8).Value = 'day'
$dayWs.Cells.Item(607,9).Value = '23/09/2024 11:34:50'

# Row 608: GMRINFRA
$dayWs.Cells.Item(608,1).Value = 25
$dayWs.Cells.Item(608,2).Value = 'GMRINFRA'
$dayWs.Cells.Item(608,3).Value = 'Gmr Infrastructure Limited'
$dayWs.Cells.Item(608,4).NumberFormat = "@"
$dayWs.Cells.Item(608,4).Value = '532754'
$dayWs.Cells.Item(608,4).ClearFormats()
$dayWs.Cells.Item(608,5).Value = 0.8
$dayWs.Cells.Item(608,6).Value = 95.29
$dayWs.Cells.Item(608,7).Value = 45233507
$dayWs.Cells.Item(608,8).Value = 'day'
$dayWs.Cells.Item(608,9).Value = '23/09/2024 11:34:50'

# --- Convert bsecode column D (rows 323-354) in the 'week' sheet from text to numeric ---

$weekWs.Cells.Item(323,4).Value = 500530
$weekWs.Cells.Item(324,4).Value = 500488
$weekWs.Cells.Item(325,4).Value = 532538
$weekWs.Cells.Item(326,4).Value = 500034
$weekWs.Cells.Item(327,4).Value = 500182
$weekWs.Cells.Item(328,4).Value = 532644
$weekWs.Cells.Item(329,4).Value = 500331
$weekWs.Cells.Item(330,4).Value = 502355
$weekWs.Cells.Item(331,4).Value = 500300
$weekWs.Cells.Item(332,4).Value = 500410
$weekWs.Cells.Item(333,4).Value = 532478
$weekWs.Cells.Item(334,4).Value = 532830
$weekWs.Cells.Item(335,4).Value = 532978
$weekWs.Cells.Item(336,4).Value = 533309
$weekWs.Cells.Item(337,4).Value = 511243
$weekWs.Cells.Item(338,4).Value = 532424
$weekWs.Cells.Item(339,4).Value = 532174
$weekWs.Cells.Item(340,4).Value = 500271
$weekWs.Cells.Item(341,4).Value = 500302
$weekWs.Cells.Item(342,4).Value = 532286
$weekWs.Cells.Item(343,4).Value = 543066
$weekWs.Cells.Item(344,4).Value = 540611
$weekWs.Cells.Item(345,4).Value = 500253
$weekWs.Cells.Item(346,4).Value = 500096
$weekWs.Cells.Item(347,4).Value = 500877
$weekWs.Cells.Item(348,4).Value = 500312
$weekWs.Cells.Item(349,4).Value = 540065
$weekWs.Cells.Item(350,4).Value = 500469
$weekWs.Cells.Item(351,4).Value = 532234
$weekWs.Cells.Item(352,4).Value = 533519
$weekWs.Cells.Item(353,4).Value = 500470
$weekWs.Cells.Item(354,4).Value = 539437
